$d = $word.ActiveDocument

# Locate the "LOM3057: Introdução aos Materiais Poliméricos (Requisito
# fraco)" paragraph (the last line of the "Requisitos" section) and
# then remove the three paragraphs that used to follow it on the
# rendered page: a blank paragraph, the "Ver no Jupiter Salvar em pdf
# Salvar em docx" paragraph, and the "© 2020 . Contact: ..." footer
# paragraph. The blank paragraph and page-break paragraph that close
# out the document stay untouched.
$rng = $d.Content
$found = $rng.Find.Execute(
    "LOM3057: Introdução aos Materiais Poliméricos (Requisito fraco)",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Find.Execute narrows $rng down to just the matched text; re-grab
    # the enclosing paragraph (by index) to get its full extent,
    # including the trailing paragraph mark.
    $reqParaIndex = $rng.Paragraphs.Item(1).Index
    $reqPara = $d.Paragraphs.Item($reqParaIndex)

    $firstToRemove = $reqPara.Next()
    $lastToRemove = $firstToRemove.Next().Next()

    $deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $deleteRange.Delete()
}
